# Update "想去人数" (want-to-go count) values in column F
# on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value  = 446
$ws1.Range("F10").Value = 7012
$ws1.Range("F16").Value = 1761
$ws1.Range("F21").Value = 116
$ws1.Range("F22").Value = 58
$ws1.Range("F26").Value = 166
$ws1.Range("F27").Value = 4113

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F12").Value = 446
$ws4.Range("F15").Value = 7012
$ws4.Range("F21").Value = 1761
$ws4.Range("F26").Value = 116
$ws4.Range("F27").Value = 58
$ws4.Range("F31").Value = 166
$ws4.Range("F32").Value = 4113
